# Update source data on the "PO List" sheet. The "Assignment" sheet's
# B10/E10 cells hold array FILTER/MAX formulas over these columns, so
# their cached results ("Jimmy Wong" -> "Chris Chong", "Brian Phua" ->
# "Michael Sien") update automatically on recalculation - we must not
# touch those cells directly, or their formulas would be overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO List")

# Row 4 (Charmaine Fang)
$ws.Range("K4").Value = 15
$ws.Range("V4").Value = 2

# Row 5 (Ng Kay Beng)
$ws.Range("K5").Value = 18

# Row 7 (Lim Kheng Guan)
$ws.Range("K7").Value = 21

# Row 8 (Vance Kang)
$ws.Range("K8").Value = 17

# Row 9 (Michael Sien)
$ws.Range("V9").Value = 6

# Row 11 (Woo Kwan Wye)
$ws.Range("V11").Value = 5

# Row 12 (Brian Phua)
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 44873
$ws.Range("K12").Value = 2
$ws.Range("S12").Value = 7
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 44873
$ws.Range("V12").Value = 1

# Row 14 (Yau Yen Nee)
$ws.Range("K14").Value = 16

# Row 16 (Tan Chong Lin)
$ws.Range("V16").Value = 4

# Row 18 (Yan Chong Hui)
$ws.Range("K18").Value = 20
$ws.Range("P18").Value = 3
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 44879

# Row 19 (Kang Mi)
$ws.Range("V19").Value = 3

# Row 20 (Bong Sell Feng)
$ws.Range("K20").Value = 19

# Row 21 (Santoso)
$ws.Range("K21").Value = 22

# Row 25 (Jimmy Wong)
$ws.Range("K25").Value = 23

# Row 27 (Chris Chong)
$ws.Range("H27").Value = 2
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 44290
$ws.Range("K27").Value = 24

$excel.CalculateFull()
